# Applies the "done git part of doc" edit:
# - Expands paragraph 3 (the one holding the _GoBack bookmark) with the
#   full narrative text, split into multiple runs around the existing
#   bookmark, with a spell-check proofErr wrap around "credential.username".
# - Appends two new empty paragraphs (matching the indent/lang formatting)
#   after that paragraph.

$d = $word.ActiveDocument

# --- Step 1: splice the new runs in around the existing _GoBack bookmark ---
# We round-trip the whole document content through WordOpenXML so we can do
# a precise, surgical string replacement while InsertXML takes care of
# re-parsing/re-applying it (this preserves every other part of the package
# untouched: styles, numbering, settings, etc).
$full = $d.Content.WordOpenXML

$oldBookmark = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$newBookmarkBlock = '<w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>I created the assignment 1 repository which is initialized with README.md file on my GitHub account and cloned it to my local, I used ‘ng new assignment1’</w:t></w:r><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> to initialized all the folder structure in my </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve">local folder and used </w:t></w:r><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve">git config </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>credential.username</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> "Wendy3813ict"</w:t></w:r><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> command to config my username and password then I have the permission to push my repository to my Git account, then I used ‘git add .’ to initialize each update commit and used ‘git commit -m ‘#comment content’’ to comment each commit and used ‘git push origin master’ to push/upload each commit, sometimes I used ‘git status’ to check my current modified files.</w:t></w:r>'

if ($full.IndexOf($oldBookmark) -lt 0) {
    throw "edit.ps1: could not locate the _GoBack bookmark anchor in Content.WordOpenXML"
}

$updated = $full.Replace($oldBookmark, $newBookmarkBlock)
$d.Content.InsertXML($updated)

# --- Step 2: append the two new empty paragraphs after that paragraph ---
$targetPara = $d.Paragraphs(3)
$endOfPara = $targetPara.Range.End
$insertionPoint = $d.Range($endOfPara, $endOfPara)

$emptyParaXml = '<w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:lang w:val="en-AU"/></w:rPr></w:pPr></w:p>'
$twoEmptyParas = $emptyParaXml + $emptyParaXml

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  $twoEmptyParas +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($packageXml)

Write-Output ("edit.ps1 done - paragraphs now: " + $d.Paragraphs.Count)
